$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# --- Crime statistics table updates (rows 15-30) ---
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -25.423728813559
$ws.Range("L16").Value = -22.807017543859
$ws.Range("M16").Value = -32.307692307692
$ws.Range("N16").Value = -86.666666666666
$ws.Range("C17").Value = 9
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 72.727272727272
$ws.Range("I17").Value = 197
$ws.Range("J17").Value = 173
$ws.Range("K17").Value = 13.872832369942
$ws.Range("L17").Value = 57.6
$ws.Range("M17").Value = 131.764705882353
$ws.Range("N17").Value = -23.046875
$ws.Range("C18").Value = 3
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 35
$ws.Range("K18").Value = -5.405405405405
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = -50.704225352112
$ws.Range("N18").Value = -91.566265060241
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = -31.25
$ws.Range("I19").Value = 86
$ws.Range("J19").Value = 86
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = -12.244897959183
$ws.Range("M19").Value = 65.384615384615
$ws.Range("N19").Value = -47.239263803681
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("I20").Value = 18
$ws.Range("K20").Value = -43.75
$ws.Range("L20").Value = -55
$ws.Range("M20").Value = -58.139534883720
$ws.Range("N20").Value = -93.050193050193
$ws.Range("C21").Value = 15
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 1.612903225806
$ws.Range("I21").Value = 402
$ws.Range("J21").Value = 402
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 9.536784741144
$ws.Range("M21").Value = 21.450151057401
$ws.Range("N21").Value = -72.294968986905
$ws.Range("L22").Value = -88.888888888888
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 46
$ws.Range("K23").Value = -8.695652173913
$ws.Range("L23").Value = -12.5
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 85.714285714285
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 53.191489361702
$ws.Range("I24").Value = 418
$ws.Range("J24").Value = 323
$ws.Range("K24").Value = 29.411764705882
$ws.Range("L24").Value = 38.410596026490
$ws.Range("M24").Value = 107.960199004975
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 160
$ws.Range("I25").Value = 66
$ws.Range("J25").Value = 26
$ws.Range("K25").Value = 153.846153846154
$ws.Range("L25").Value = 22.222222222222
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 31.578947368421
$ws.Range("I26").Value = 265
$ws.Range("J26").Value = 220
$ws.Range("K26").Value = 20.454545454545
$ws.Range("L26").Value = 15.720524017467
$ws.Range("M26").Value = 19.909502262443
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 4.545454545454
$ws.Range("L27").Value = 21.052631578947
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 22
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -4.347826086956
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 8
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 14.285714285714
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -61.904761904761
$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 6
$ws.Range("K30").Value = -14.285714285714
$ws.Range("L30").Value = 20
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -71.428571428571
Write-Output "edit complete"
